$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 with new sensor readings (columns A:AH)
$row2 = @(45176.50694444445,20.658,14.028,4.08,43.704,35.637,16.257,52.557,25.014,10.522,16.084,17.263,18.024,5.19,16.166,22.644,13.689,3.709,2.453,238.72,44.979,14.922,29.658,15.316,2.936,26.039,13.181,11.919,13.951,17.711,3.458,46.546,8.337999999999999,18.656)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @(45176.51388888889,23.541,16.982,2.047,50.805,41.68,18.525,70.878,28.505,12.477,18.64,20.381,21.369,5.918,18.422,26.097,15.636,1.663,1.291,273.097,51.526,17.004,34.404,18.101,2.921,34.387,15.02,13.45,15.781,21.162,1.266,64.35299999999999,9.561999999999999,21.259)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @(45176.52083333334,18.737,13.617,1.399,40.532,33.245,14.745,58.718,22.687,9.981999999999999,14.844,16.297,17.09,4.711,14.663,20.801,12.453,1.092,0.89,215.869,41.048,13.534,27.445,14.457,2.288,28.077,11.955,10.684,12.54,16.96,0.773,53.274,7.606,16.921)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @(45176.52777777778,3.84,2.55,0.68,8.210000000000001,6.56,3.03,18,4.65,2,2.81,3.33,3.4,0.98,3.01,4.27,2.75,0.72,0.34,38.53,8.77,2.78,5.73,2.96,0.68,7.97,2.45,2.31,2.69,3.38,0.5600000000000001,16.75,1.49,3.47)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# Remove old row 6 (data now only spans rows 1-5)
$ws.Rows.Item(6).Delete()

# Widen columns J (10) and X (24) from 7 to 8 characters
$ws.Columns.Item(10).ColumnWidth = 7.14
$ws.Columns.Item(24).ColumnWidth = 7.14

